$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to move from 45207 (2023-10-08)
# to 45208 (2023-10-09) for rows 2 through 10.
foreach ($row in 2..10) {
    $ws.Cells.Item($row, 3).Value = 45208
}
